$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Polls" sheet: append two new poll rows (17 and 18) and center-align the
#    numeric/result columns (D:P) for the whole data body (existing + new).
# ---------------------------------------------------------------------------
$polls = $wb.Worksheets.Item("Polls")

# Row 17 - MassiveCaller poll (2022-05-03)
$polls.Cells.Item(17, 1).Value = "MassiveCaller"
$polls.Cells.Item(17, 2).Value = 44684
$polls.Cells.Item(17, 2).Style = "Comma"
$polls.Cells.Item(17, 2).NumberFormat = "yyyy\-mm\-dd;@"
$polls.Cells.Item(17, 3).Value = 1000
$polls.Cells.Item(17, 4).Value = 0.034
$polls.Cells.Item(17, 5).Value = 0.352
$polls.Cells.Item(17, 6).Value = 0.304
$polls.Cells.Item(17, 7).Value = 0.081
$polls.Cells.Item(17, 8).Value = 0.088
$polls.Cells.Item(17, 9).Value = 0.015
$polls.Cells.Item(17, 10).Value = 0.01
$polls.Cells.Item(17, 11).Value = 0.01
$polls.Cells.Item(17, 12).Value = "NA"
$polls.Cells.Item(17, 13).Value = 0.142
$polls.Cells.Item(17, 14).Value = "NA"
$polls.Cells.Item(17, 15).Value = "NA"
$polls.Cells.Item(17, 16).Value = "MassiveCaller"
$polls.Cells.Item(17, 17).Value = "https://www.massivecaller.com/files/colombia.pdf"
$polls.Hyperlinks.Add($polls.Cells.Item(17, 17), "https://www.massivecaller.com/files/colombia.pdf")
$polls.Cells.Item(17, 17).Style = "Hyperlink"

# Row 18 - Yanhaas poll (2022-05-04)
$polls.Cells.Item(18, 1).Value = "Yanhaas"
$polls.Cells.Item(18, 2).Value = 44685
$polls.Cells.Item(18, 2).Style = "Comma"
$polls.Cells.Item(18, 2).NumberFormat = "yyyy\-mm\-dd;@"
$polls.Cells.Item(18, 3).Value = 1232
$polls.Cells.Item(18, 4).Value = 0.032
$polls.Cells.Item(18, 5).Value = 0.4
$polls.Cells.Item(18, 6).Value = 0.21
$polls.Cells.Item(18, 7).Value = 0.07
$polls.Cells.Item(18, 8).Value = 0.12
$polls.Cells.Item(18, 9).Value = 0.01
$polls.Cells.Item(18, 10).Value = "NA"
$polls.Cells.Item(18, 11).Value = 0.003
$polls.Cells.Item(18, 12).Value = 0.003
$polls.Cells.Item(18, 13).Value = 0.13
$polls.Cells.Item(18, 14).Value = "NA"
$polls.Cells.Item(18, 15).Value = 0.06
$polls.Cells.Item(18, 16).Value = "Guarumo"
$polls.Cells.Item(18, 17).Value = "https://www.eltiempo.com/elecciones-2022/presidencia/encuesta-presidencial-petro-se-mantiene-primero-fico-segundo-671401"
$polls.Hyperlinks.Add($polls.Cells.Item(18, 17), "https://www.eltiempo.com/elecciones-2022/presidencia/encuesta-presidencial-petro-se-mantiene-primero-fico-segundo-671401")
$polls.Cells.Item(18, 17).Style = "Hyperlink"

# Center-align columns D:P across the whole data body (rows 2-18), matching
# the formatting pass applied when the new rows were added.
$polls.Range("D2:P18").HorizontalAlignment = -4108

# Selection ends up on Q19 (just past the new last row) and Polls is no
# longer the active tab.
$polls.Range("Q19").Select()

# ---------------------------------------------------------------------------
# 2. "Info" sheet: refresh the source link (new La Silla Vacía article) and
#    make it an active hyperlink; this becomes the active sheet on save.
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Item("Info")
$newLink = "https://www.lasillavacia.com/historias/silla-nacional/el-semaforo-de-las-encuestadoras-actualizado-de-cara-a-primera-vuelta/"
$info.Range("B2").Value = $newLink
$info.Hyperlinks.Add($info.Range("B2"), $newLink)
$info.Range("B2").Style = "Hyperlink"
$info.Columns.Item(2).AutoFit()

$info.Activate()
$info.Range("B11").Select()
